$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of rows 2<->4 and rows 3<->5 across columns
# D, I, J, K, L, M, N, P, Q (the columns whose values actually differ
# between the two rows in each pair).
$cols = @("D", "I", "J", "K", "L", "M", "N", "P", "Q")

function Swap-Rows($ws, $rowA, $rowB, $cols) {
    foreach ($col in $cols) {
        $rangeA = $ws.Range("$col$rowA")
        $rangeB = $ws.Range("$col$rowB")
        $valA = $rangeA.Value2
        $valB = $rangeB.Value2
        $rangeA.Value2 = $valB
        $rangeB.Value2 = $valA
    }
}

Swap-Rows $ws 2 4 $cols
Swap-Rows $ws 3 5 $cols
